$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "Voluntary acknowledgment of paternity - VAP"
$ws.Range("B18").Value = "https://www.illinoislegalaid.org/legal-information/voluntary-acknowledgment-parentage-vap"
$ws.Hyperlinks.Add($ws.Range("B18"), "https://www.illinoislegalaid.org/legal-information/voluntary-acknowledgment-parentage-vap")
$ws.Range("B18").Style = "Hyperlink"
